$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Scheduled Time" (F) values for rows 14-21 ---
$ws.Range("F14").Value = 45328.832638888889
$ws.Range("F15").Value = 45336.568749999999
$ws.Range("F16").Value = 45336.568749999999
$ws.Range("F17").Value = 45336.568749999999
$ws.Range("F18").Value = 45336.568749999999
$ws.Range("F19").Value = 45336.568749999999
$ws.Range("F20").Value = 45336.568749999999
$ws.Range("F21").Value = 45336.652083333334

# --- Append new rows 22-28 (Post IDs 21-27), Facebook "TestFI" entries ---
$newRows = @(
    @{ Row=22; Id=21; Content="TestFI 12"; Img="img15" },
    @{ Row=23; Id=22; Content="TestFI 13"; Img="img16" },
    @{ Row=24; Id=23; Content="TestFI 14"; Img="img17" },
    @{ Row=25; Id=24; Content="TestFI 15"; Img="img18" },
    @{ Row=26; Id=25; Content="TestFI 16"; Img="img19" },
    @{ Row=27; Id=26; Content="TestFI 17"; Img="img20" },
    @{ Row=28; Id=27; Content="TestFI 18"; Img="img21" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Id
    $ws.Range("B$row").Value = "Facebook"
    $ws.Range("C$row").Value = $r.Content
    $ws.Range("D$row").Value = $r.Img
    $ws.Range("E$row").Value = "#new #tech #insta"
    $ws.Range("F$row").Value = 45336.652083333334
    $ws.Range("F$row").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
    $ws.Range("G$row").Value = "Scheduled"
}

# --- Update sheet view: clear frozen/scrolled topLeftCell, update selection ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F21:F28").Select()
